$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 76926440
$ws.Range("I76").Value = 3750.5
$ws.Range("J76").Value = 142860180
$ws.Range("K76").Value = 3750.5
$ws.Range("L76").Value = 142860180
$ws.Range("M76").Value = -3435.5
$ws.Range("N76").Value = -142860810

$ws.Range("H79").Value = 76926440
$ws.Range("I79").Value = 3750.5
$ws.Range("J79").Value = 142860180
$ws.Range("K79").Value = 3750.5
$ws.Range("L79").Value = 142860180
$ws.Range("M79").Value = -2658.5
$ws.Range("N79").Value = -142862364

$ws.Range("H88").Value = 1031029.8
$ws.Range("I88").Value = 890
$ws.Range("J88").Value = 1237057.8
$ws.Range("K88").Value = 890
$ws.Range("L88").Value = 1237057.8
$ws.Range("M88").Value = -484
$ws.Range("N88").Value = -1237869.8

$ws.Range("H91").Value = 1031029.8
$ws.Range("I91").Value = 890
$ws.Range("J91").Value = 1237057.8
$ws.Range("K91").Value = 890
$ws.Range("L91").Value = 1237057.8
$ws.Range("M91").Value = 514
$ws.Range("N91").Value = -1239865.8

$ws.Range("H107").Value = 2230.1177
$ws.Range("I107").Value = 1700.5834
$ws.Range("K107").Value = 1700.5834
$ws.Range("M107").Value = 219.4166

$ws.Range("H132").Value = 6294563.5
$ws.Range("I132").Value = 9526093
$ws.Range("J132").Value = 11034.277
$ws.Range("K132").Value = 28578279
$ws.Range("L132").Value = 33102.831
$ws.Range("M132").Value = -28575749
$ws.Range("N132").Value = -38162.831

$ws.Range("H138").Value = 822174
$ws.Range("I138").Value = 1235.8667
$ws.Range("J138").Value = 1146228.5
$ws.Range("K138").Value = 3707.6001
$ws.Range("L138").Value = 3438685.5
$ws.Range("M138").Value = 1432.3999
$ws.Range("N138").Value = -3448965.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13698
$ws.Range("I2").Value = 1109
$ws.Range("J2").Value = 34679.668
$ws.Range("K2").Value = 1109
$ws.Range("L2").Value = 34679.668
$ws.Range("M2").Value = -996
$ws.Range("N2").Value = -34905.668

$ws.Range("H32").Value = 4994.926
$ws.Range("I32").Value = 5248.7085
$ws.Range("K32").Value = 5248.7085
$ws.Range("M32").Value = -4961.7085

$ws.Range("H74").Value = 1756
$ws.Range("I74").Value = 1341.3334
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 1341.3334
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -467.3334
$ws.Range("N74").Value = -4748

$ws.Range("H77").Value = 1756
$ws.Range("I77").Value = 1341.3334
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 6706.666999999999
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -2338.666999999999
$ws.Range("N77").Value = -23736

$ws.Range("H116").Value = 13698
$ws.Range("I116").Value = 1109
$ws.Range("J116").Value = 34679.668
$ws.Range("K116").Value = 1109
$ws.Range("L116").Value = 34679.668
$ws.Range("M116").Value = 1185
$ws.Range("N116").Value = -39267.668

$ws.Range("H122").Value = 1116.9
$ws.Range("I122").Value = 1116.9
$ws.Range("K122").Value = 3350.7
$ws.Range("M122").Value = -900.7000000000003


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13698
$ws.Range("I3").Value = 1109
$ws.Range("J3").Value = 34679.668
$ws.Range("K3").Value = 1109
$ws.Range("L3").Value = 34679.668
$ws.Range("M3").Value = -995
$ws.Range("N3").Value = -34907.668

$ws.Range("H20").Value = 1805.4546
$ws.Range("I20").Value = 1842.5385
$ws.Range("K20").Value = 1842.5385
$ws.Range("M20").Value = -1595.5385

$ws.Range("H107").Value = 1386.7
$ws.Range("I107").Value = 982.82355
$ws.Range("J107").Value = 1914.8462
$ws.Range("K107").Value = 982.82355
$ws.Range("L107").Value = 1914.8462
$ws.Range("M107").Value = 937.17645
$ws.Range("N107").Value = -5754.8462

$ws.Range("H134").Value = 7124.3887
$ws.Range("I134").Value = 1349.6364
$ws.Range("J134").Value = 16199
$ws.Range("K134").Value = 4048.9092
$ws.Range("L134").Value = 48597
$ws.Range("M134").Value = -1513.9092
$ws.Range("N134").Value = -53667


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71429490
$ws.Range("I16").Value = 83334190
$ws.Range("K16").Value = 83334190
$ws.Range("M16").Value = -83333903

$ws.Range("H31").Value = 1705.8462
$ws.Range("I31").Value = 1705.8462
$ws.Range("K31").Value = 1705.8462
$ws.Range("M31").Value = -1410.8462

$ws.Range("H34").Value = 1705.8462
$ws.Range("I34").Value = 1705.8462
$ws.Range("K34").Value = 1705.8462
$ws.Range("M34").Value = -1503.8462

$ws.Range("H113").Value = 71429490
$ws.Range("I113").Value = 83334190
$ws.Range("K113").Value = 83334190
$ws.Range("M113").Value = -83332020

$ws.Range("H132").Value = 5126.5713
$ws.Range("I132").Value = 5324.7915
$ws.Range("J132").Value = 3937.25
$ws.Range("K132").Value = 15974.3745
$ws.Range("L132").Value = 11811.75
$ws.Range("M132").Value = -13444.3745
$ws.Range("N132").Value = -16871.75


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 576.5333000000001
$ws.Range("I107").Value = 346.5
$ws.Range("J107").Value = 660.1818
$ws.Range("K107").Value = 1039.5
$ws.Range("L107").Value = 1980.5454
$ws.Range("M107").Value = 880.5
$ws.Range("N107").Value = -5820.5454

$ws.Range("H131").Value = 911.4
$ws.Range("J131").Value = 956.01086
$ws.Range("L131").Value = 2868.03258
$ws.Range("N131").Value = -12948.03258

$ws.Range("H132").Value = 1059.8334
$ws.Range("J132").Value = 1333
$ws.Range("L132").Value = 11997
$ws.Range("N132").Value = -17057


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2169.6667
$ws.Range("I102").Value = 2312.2
$ws.Range("J102").Value = 1457
$ws.Range("K102").Value = 2312.2
$ws.Range("L102").Value = 1457
$ws.Range("M102").Value = -690.1999999999998
$ws.Range("N102").Value = -4701

$ws.Range("H132").Value = 2907.6667
$ws.Range("I132").Value = 3125.375
$ws.Range("J132").Value = 2472.25
$ws.Range("K132").Value = 9376.125
$ws.Range("L132").Value = 7416.75
$ws.Range("M132").Value = -6846.125
$ws.Range("N132").Value = -12476.75

$ws.Range("H136").Value = 16500
$ws.Range("J136").Value = 16500
$ws.Range("L136").Value = 49500
$ws.Range("N136").Value = -54600


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1328.9524
$ws.Range("I61").Value = 1116.6428
$ws.Range("J61").Value = 1753.5714
$ws.Range("K61").Value = 1116.6428
$ws.Range("L61").Value = 1753.5714
$ws.Range("M61").Value = -914.6428000000001
$ws.Range("N61").Value = -2157.5714

$ws.Range("H113").Value = 1328.9524
$ws.Range("I113").Value = 1116.6428
$ws.Range("J113").Value = 1753.5714
$ws.Range("K113").Value = 1116.6428
$ws.Range("L113").Value = 1753.5714
$ws.Range("M113").Value = 1053.3572
$ws.Range("N113").Value = -6093.5714

$ws.Range("H122").Value = 14708974
$ws.Range("I122").Value = 25003184
$ws.Range("J122").Value = 2958.2856
$ws.Range("K122").Value = 75009552
$ws.Range("L122").Value = 8874.856800000001
$ws.Range("M122").Value = -75007102
$ws.Range("N122").Value = -13774.8568

$ws.Range("H132").Value = 80064.766
$ws.Range("J132").Value = 127981.625
$ws.Range("L132").Value = 383944.875
$ws.Range("N132").Value = -389004.875


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 40005.5
$ws.Range("J20").Value = 40005.5
$ws.Range("L20").Value = 40005.5
$ws.Range("N20").Value = -40485.5

$ws.Range("H113").Value = 1001
$ws.Range("I113").Value = 767.3333
$ws.Range("J113").Value = 1351.5
$ws.Range("K113").Value = 2301.9999
$ws.Range("L113").Value = 4054.5
$ws.Range("M113").Value = -131.9998999999998
$ws.Range("N113").Value = -8394.5

$ws.Range("H123").Value = 57857
$ws.Range("J123").Value = 57857
$ws.Range("N123").Value = -67657


Write-Output "done"